$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 3-5 (B column) with revised inflation rates
$ws.Range("B3").Value = 0.01035
$ws.Range("B4").Value = 0.010038
$ws.Range("B5").Value = 0.01025

# Row 6: TimeToZero changes from 10 to 4, rate updated
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 0.01046

# New rows 7-13 with additional curve points
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 0.010691

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 0.011175

$ws.Range("A9").Value = 10
$ws.Range("B9").Value = 0.012

$ws.Range("A10").Value = 15
$ws.Range("B10").Value = 0.01316

$ws.Range("A11").Value = 20
$ws.Range("B11").Value = 0.01411

$ws.Range("A12").Value = 30
$ws.Range("B12").Value = 0.015387

$ws.Range("A13").Value = 50
$ws.Range("B13").Value = 0.015825

# Match the percentage number format used by the rest of column B
$ws.Range("B7:B13").NumberFormat = "0.00%"

# Update the active selection to reflect the new last-used row
$ws.Range("B14").Select() | Out-Null
